$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.692.48"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.864.96"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("D4").Value = "'1.012"
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").Value = "'334.13"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").Value = "'1.010"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").Value = "'0.4703"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.3924"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "'45.92"
$ws.Range("E9").Value = "  -3.86%  "
$ws.Range("D10").Value = "'0.08003"
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("D11").Value = "'1.005"
$ws.Range("E11").Value = "  -2.20%  "
$ws.Range("D12").Value = "'21.82"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("D13").Value = "1.875.67"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").Value = "'5.999"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "'7.259"
$ws.Range("E15").Value = "  +1.80%  "
$ws.Range("D16").Value = "'1.012"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "'88.41"
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").Value = "'0.06729"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "'0.00001045"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "'17.05"
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("D21").Value = "'1.011"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "27.662.59"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").Value = "'5.470"
$ws.Range("E23").Value = "  -1.24%  "
$ws.Range("D24").Value = "'10.91"
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("D25").Value = "'2.309"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").Value = "2.102.50"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").Value = "'159.50"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "'19.83"
$ws.Range("E28").Value = "  -1.97%  "
$ws.Range("D29").Value = "'2.145"
$ws.Range("E29").Value = "  +1.99%  "
$ws.Range("D30").Value = "'5.446"
$ws.Range("E30").Value = "  -2.44%  "
$ws.Range("D31").Value = "'121.75"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "'0.9821"
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("D33").Value = "'0.09504"
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").Value = "'3.615"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").Value = "'5.303"
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("D36").Value = "'1.336"
$ws.Range("E36").Value = "  -8.12%  "
$ws.Range("D37").Value = "'0.06068"
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("D38").Value = "'0.02232"
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("D39").Value = "'8.321"
$ws.Range("E39").Value = "  +2.44%  "
$ws.Range("D40").Value = "'1.196"
$ws.Range("E40").Value = "  -2.93%  "
$ws.Range("D41").Value = "'1.010"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").Value = "'0.5966"
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("D43").Value = "'0.1887"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D44").Value = "'10.29"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("D45").Value = "'1.252"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("D46").Value = "'0.5647"
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("D47").Value = "'12.19"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").Value = "'1.927"
$ws.Range("E48").Value = "  -1.15%  "
$ws.Range("D49").Value = "'3.190"
$ws.Range("E49").Value = "  -6.11%  "
$ws.Range("D50").Value = "'0.06769"
$ws.Range("E50").Value = "  -2.07%  "
$ws.Range("D51").Value = "'112.21"
$ws.Range("E51").Value = "  -2.13%  "
